# "minor corrections to several slides"
#
# Slide 58 ("Method emit() for Class Variable") - Content Placeholder 2:
# the paragraph describing the method signature currently reads
#   "public void emit() throws CodeGenException, IOException"
# and should be corrected to drop the erroneous ", IOException" suffix:
#   "public void emit() throws CodeGenException"

$p = $ppt.ActivePresentation

$needle = "public void emit() throws CodeGenException, IOException"
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*$needle*") {
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange

    # Find the paragraph that contains the method signature. Note:
    # TextRange.Text for a single paragraph includes the trailing
    # paragraph-mark character, so trim that before comparing.
    $paraCount = $tr.Paragraphs().Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $para = $tr.Paragraphs($pi, 1)
        if ($para.Text.TrimEnd([char]13, [char]10, [char]11) -eq $needle) {
            [void]$para.Replace(", IOException", "")
            break
        }
    }
}
